$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 1 de Octubre de 2020 a las 12:38"

# 2. Row 4 - Estados Unidos: refreshed daily figures
$ws.Range("B4").Value = 7450637
$ws.Range("C4").Value = 3355
$ws.Range("E4").Value = 2538113
$ws.Range("G4").Value = 38
$ws.Range("H4").Value = 211778

# 3. Row 18 - Banglades: refreshed daily figures
$ws.Range("B18").Value = 364987
$ws.Range("C18").Value = 1508
$ws.Range("D18").Value = 277078
$ws.Range("E18").Value = 82637
$ws.Range("G18").Value = 21
$ws.Range("H18").Value = 5272

# 4. Row 32 - Rumania: refreshed daily figures
$ws.Range("B32").Value = 129658
$ws.Range("C32").Value = 2086
$ws.Range("D32").Value = 103994
$ws.Range("E32").Value = 20802
$ws.Range("G32").Value = 37
$ws.Range("H32").Value = 4862

# 5. Row 61 - Suiza: refreshed daily figures
$ws.Range("B61").Value = 53832
$ws.Range("C61").Value = 550
$ws.Range("E61").Value = 9058

# 6. Malasia overtakes Namibia in the ranking (rows 98-99 swap, Malasia gets new numbers)
$ws.Range("A98").Value = "Malasia"
$ws.Range("B98").Value = 11484
$ws.Range("C98").Value = 260
$ws.Range("D98").Value = 10014
$ws.Range("E98").Value = 1334
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 136

$ws.Range("A99").Value = "Namibia"
$ws.Range("B99").Value = 11265
$ws.Range("C99").Value = 0
$ws.Range("D99").Value = 9014
$ws.Range("E99").Value = 2130
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 121

# 7. Nueva Caledonia swaps ahead of Santa Lucia (rows 207-208, tied totals, figures unchanged)
$ws.Range("A207").Value = "Nueva Caledonia"
$ws.Range("B207").Value = 27
$ws.Range("C207").Value = 0
$ws.Range("D207").Value = 27
$ws.Range("E207").Value = 0
$ws.Range("F207").Value = 0
$ws.Range("G207").Value = 0
$ws.Range("H207").Value = 0

$ws.Range("A208").Value = "Santa Lucia"
$ws.Range("B208").Value = 27
$ws.Range("C208").Value = 0
$ws.Range("D208").Value = 27
$ws.Range("E208").Value = 0
$ws.Range("F208").Value = 0
$ws.Range("G208").Value = 0
$ws.Range("H208").Value = 0
